# Auto-generated script applying the Typhon_Profits market-data refresh diff
$wb = $excel.ActiveWorkbook

function Set-Cells($ws, $values) {
    foreach ($ref in $values.Keys) {
        $v = $values[$ref]
        if ($null -eq $v) {
            $ws.Range($ref).ClearContents()
        } else {
            $ws.Range($ref).Value = $v
        }
    }
}

$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H13" = 0
    "J13" = 0
    "L13" = $null
    "N13" = 0
    "H28" = 413.63635
    "I28" = 162.25
    "J28" = 1084
    "K28" = 162.25
    "L28" = 1084
    "M28" = 322.75
    "N28" = -2054
    "H96" = 19231088
    "I96" = 19231088
    "K96" = 57693264
    "M96" = -57691891
    "H97" = 1115.7142
    "J97" = 1115.7142
    "L97" = 3347.1426
    "N97" = -4339.142599999999
    "H98" = 717.9032
    "I98" = 768.05
    "K98" = 768.05
    "M98" = 729.95
    "H112" = 1129.1538
    "J112" = 1129.1538
    "L112" = 3387.4614
    "N112" = -5603.4614
    "H122" = 717.9032
    "I122" = 768.05
    "K122" = 2304.15
    "M122" = 145.8500000000004
    "H125" = 1163.375
    "J125" = 1163.375
    "L125" = 10470.375
    "N125" = -15390.375
    "H137" = 30830.412
    "I137" = 1418.1364
    "J137" = 84752.914
    "K137" = 4254.4092
    "L137" = 254258.742
    "M137" = -1704.4092
    "N137" = -259358.742
    "H138" = 2000.5641
    "J138" = 2219.2769
    "L138" = 6657.8307
    "N138" = -16937.8307
}
Set-Cells $ws $values

$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H32" = 26945.936
    "I32" = 33305.297
    "J32" = 3416.3
    "K32" = 33305.297
    "L32" = 3416.3
    "M32" = -33018.297
    "N32" = -3990.3
    "H61" = 2354.9656
    "I61" = 2088.8096
    "K61" = 2088.8096
    "M61" = -1876.8096
    "H63" = 3475330.5
    "I63" = 1696
    "J63" = 5212148
    "K63" = 1696
    "L63" = 5212148
    "M63" = -1010
    "N63" = -5213520
    "H66" = 3475330.5
    "I66" = 1696
    "J66" = 5212148
    "K66" = 8480
    "L66" = 26060740
    "M66" = -5048
    "N66" = -26067604
    "H74" = 40001148
    "I74" = 58824144
    "J74" = 2278.125
    "K74" = 58824144
    "L74" = 2278.125
    "M74" = -58823270
    "N74" = -4026.125
    "H77" = 40001148
    "I77" = 58824144
    "J77" = 2278.125
    "K77" = 294120720
    "L77" = 11390.625
    "M77" = -294116352
    "N77" = -20126.625
    "H97" = 861.35895
    "I97" = 850.5484
    "K97" = 850.5484
    "M97" = -354.5484
    "H102" = 1344.9
    "I102" = 1074.8334
    "K102" = 1074.8334
    "M102" = 547.1666
    "H110" = 578.625
    "I110" = 365
    "J110" = 649.8333
    "K110" = 365
    "L110" = 649.8333
    "M110" = 1680
    "N110" = -4739.8333
    "H122" = 3145.4443
    "I122" = 2044.2858
    "J122" = 6999.5
    "K122" = 6132.857400000001
    "L122" = 20998.5
    "M122" = -3682.857400000001
    "N122" = -25898.5
    "H136" = 2354.9656
    "I136" = 2088.8096
    "K136" = 6266.4288
    "M136" = -3716.4288
}
Set-Cells $ws $values

$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H99" = 1781.619
    "I99" = 1533.5333
    "J99" = 2401.8333
    "K99" = 1533.5333
    "L99" = 2401.8333
    "M99" = -35.53330000000005
    "N99" = -5397.8333
    "H105" = 3128665.5
    "I105" = 4004.9
    "J105" = 8336433.5
    "K105" = 4004.9
    "L105" = 8336433.5
    "M105" = -2257.9
    "N105" = -8339927.5
    "H107" = 1693.0667
    "I107" = 1398.6666
    "J107" = 2870.6667
    "K107" = 1398.6666
    "L107" = 2870.6667
    "M107" = 521.3334
    "N107" = -6710.6667
    "H134" = 34052.637
    "I134" = 42836.117
    "J134" = 1428.2858
    "K134" = 128508.351
    "L134" = 4284.857400000001
    "M134" = -125973.351
    "N134" = -9354.857400000001
}
Set-Cells $ws $values

$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H16" = 499.77777
    "I16" = 499.7143
    "K16" = 499.7143
    "M16" = -212.7143
    "H31" = 11619.229
    "I31" = 28168.666
    "K31" = 28168.666
    "M31" = -27873.666
    "H34" = 11619.229
    "I34" = 28168.666
    "K34" = 28168.666
    "M34" = -27966.666
    "H105" = 25000850
    "I105" = 25000850
    "K105" = 25000850
    "M105" = -24999103
    "H110" = 40998.75
    "J110" = 40998.75
    "L110" = 40998.75
    "N110" = -49178.75
    "H113" = 499.77777
    "I113" = 499.7143
    "K113" = 499.7143
    "M113" = 1670.2857
}
Set-Cells $ws $values

$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H23" = 1344.0834
    "I23" = 726.6667
    "J23" = 1549.8889
    "K23" = 2180.0001
    "L23" = 4649.6667
    "M23" = -1945.0001
    "N23" = -5119.6667
    "H107" = 5247.7
    "I107" = 16833.334
    "J107" = 282.42856
    "K107" = 50500.00199999999
    "L107" = 847.28568
    "M107" = -48580.00199999999
    "N107" = -4687.28568
    "H122" = 673.8333
    "I122" = 324.5
    "J122" = 923.3570999999999
    "K122" = 2920.5
    "L122" = 8310.213899999999
    "M122" = -470.5
    "N122" = -13210.2139
    "H131" = 789.97
    "J131" = 814.76044
    "L131" = 2444.28132
    "N131" = -12524.28132
}
Set-Cells $ws $values

$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H4" = 10000
    "J4" = 10000
    "L4" = 10000
    "N4" = -10224
    "H11" = 3326250
    "I11" = 3461538.5
    "J11" = 2740000
    "K11" = 3461538.5
    "L11" = 2740000
    "M11" = -3461399.5
    "N11" = -2740278
    "H97" = 533.56525
    "J97" = 199
    "L97" = 199
    "N97" = -1191
    "H107" = 7692639
    "I107" = 390
    "J107" = 25641220
    "K107" = 390
    "L107" = 25641220
    "M107" = 1530
    "N107" = -25645060
    "H113" = 3121.9707
    "I113" = 2883.5454
    "J113" = 3559.0833
    "K113" = 2883.5454
    "L113" = 3559.0833
    "M113" = -713.5454
    "N113" = -7899.0833
    "H132" = 63703.72
    "I132" = 71439.8
    "J132" = 52099.6
    "K132" = 214319.4
    "L132" = 156298.8
    "M132" = -211789.4
    "N132" = -161358.8
    "H136" = 15914.923
    "J136" = 15914.923
    "L136" = 47744.769
    "N136" = -52844.769
}
Set-Cells $ws $values

$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H61" = 4948.5415
    "I61" = 3217.9285
    "J61" = 7371.4
    "K61" = 3217.9285
    "L61" = 7371.4
    "M61" = -3015.9285
    "N61" = -7775.4
    "H93" = 876.4375
    "I93" = 876.4375
    "K93" = 876.4375
    "M93" = 371.5625
    "H113" = 4948.5415
    "I113" = 3217.9285
    "J113" = 7371.4
    "K113" = 3217.9285
    "L113" = 7371.4
    "M113" = -1047.9285
    "N113" = -11711.4
    "H132" = 2807.8462
    "I132" = 2280.8
    "J132" = 3137.25
    "K132" = 6842.400000000001
    "L132" = 9411.75
    "M132" = -4312.400000000001
    "N132" = -14471.75
    "H136" = 14493.111
    "I136" = 16134.125
    "K136" = 48402.375
    "M136" = -45852.375
}
Set-Cells $ws $values
